$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = 8.894566086880333
$ws.Range("C2").Value = 4.257033412227463
$ws.Range("D2").Value = 10.2410816369089
$ws.Range("F2").Value = 34.52407848195887
$ws.Range("G2").Value = 3.665015584202522
$ws.Range("I2").Value = 23.63849788701656
$ws.Range("J2").Value = 11.4373519350174
$ws.Range("K2").Value = 9.134507244806032
$ws.Range("M2").Value = 15.96144507318415
$ws.Range("O2").Value = 25.60824342392698
$ws.Range("B3").Value = 8.616210680144381
$ws.Range("C3").Value = 4.067635136018182
$ws.Range("D3").Value = 10.19494186520907
$ws.Range("F3").Value = 34.60141949757784
$ws.Range("G3").Value = 3.666778967141621
$ws.Range("I3").Value = 23.73639263940846
$ws.Range("J3").Value = 11.45709493299023
$ws.Range("K3").Value = 8.954745468390767
$ws.Range("M3").Value = 15.88167870352883
$ws.Range("O3").Value = 25.69668982894278
$ws.Range("B4").Value = 8.441786187800343
$ws.Range("C4").Value = 3.947034035510015
$ws.Range("D4").Value = 10.16828796351986
$ws.Range("F4").Value = 34.65656064225659
$ws.Range("G4").Value = 3.667919343004601
$ws.Range("I4").Value = 23.80068773081309
$ws.Range("J4").Value = 11.4708659654549
$ws.Range("K4").Value = 8.843709814135901
$ws.Range("M4").Value = 15.83503682105965
$ws.Range("O4").Value = 25.75597084787127
$ws.Range("B5").Value = 8.369932457087623
$ws.Range("C5").Value = 3.896869630099502
$ws.Range("D5").Value = 10.15785602851997
$ws.Range("F5").Value = 34.68095209816802
$ws.Range("G5").Value = 3.668398597183369
$ws.Range("I5").Value = 23.82794128628838
$ws.Range("J5").Value = 11.47689242330085
$ws.Range("K5").Value = 8.798354767817534
$ws.Range("M5").Value = 15.81663147101093
$ws.Range("O5").Value = 25.78137758531409
$ws.Range("B6").Value = 8.357957782910299
$ws.Range("C6").Value = 3.888480443538578
$ws.Range("D6").Value = 10.15615000080562
$ws.Range("F6").Value = 34.68511819304872
$ws.Range("G6").Value = 3.668479056578996
$ws.Range("I6").Value = 23.83253030015573
$ws.Range("J6").Value = 11.47791815680243
$ws.Range("K6").Value = 8.790818885899448
$ws.Range("M6").Value = 15.81361201860809
$ws.Range("O6").Value = 25.78567175673701
$ws.Range("B7").Value = 8.440820128248559
$ws.Range("C7").Value = 3.946361531594499
$ws.Range("D7").Value = 10.16814552422253
$ws.Range("F7").Value = 34.65688182091532
$ws.Range("G7").Value = 3.667925747455681
$ws.Range("I7").Value = 23.80105101930003
$ws.Range("J7").Value = 11.47094556142237
$ws.Range("K7").Value = 8.843098494546927
$ws.Range("M7").Value = 15.83478614572186
$ws.Range("O7").Value = 25.75630843611603
$ws.Range("B8").Value = 8.799381137216676
$ws.Range("C8").Value = 4.192659755303463
$ws.Range("D8").Value = 10.22482988833181
$ws.Range("F8").Value = 34.54915518342209
$ws.Range("G8").Value = 3.665611660055741
$ws.Range("I8").Value = 23.67138282629192
$ws.Range("J8").Value = 11.44381716826633
$ws.Range("K8").Value = 9.072697594767815
$ws.Range("M8").Value = 15.93346589087407
$ws.Range("O8").Value = 25.63770629058343
$ws.Range("B9").Value = 9.47028056230644
$ws.Range("C9").Value = 4.63902663906244
$ws.Range("D9").Value = 10.34890054960911
$ws.Range("F9").Value = 34.39877333221605
$ws.Range("G9").Value = 3.661529135971603
$ws.Range("I9").Value = 23.45033503137247
$ws.Range("J9").Value = 11.4036971894992
$ws.Range("K9").Value = 9.515185865979197
$ws.Range("M9").Value = 16.14483652623569
$ws.Range("O9").Value = 25.44467152478556
$ws.Range("B10").Value = 9.938299190709804
$ws.Range("C10").Value = 4.941941359610999
$ws.Range("D10").Value = 10.44737884874842
$ws.Range("F10").Value = 34.32556891277305
$ws.Range("G10").Value = 3.65880446077154
$ws.Range("I10").Value = 23.30819880466767
$ws.Range("J10").Value = 11.3821883091107
$ws.Range("K10").Value = 9.832347427884782
$ws.Range("M10").Value = 16.31008123220781
$ws.Range("O10").Value = 25.32705326384215
$ws.Range("B11").Value = 10.14486907740386
$ws.Range("C11").Value = 5.073873720653293
$ws.Range("D11").Value = 10.49364368079906
$ws.Range("F11").Value = 34.30038946676061
$ws.Range("G11").Value = 3.65762398666008
$ws.Range("I11").Value = 23.24794094294406
$ws.Range("J11").Value = 11.37413126647214
$ws.Range("K11").Value = 9.974280982269413
$ws.Range("M11").Value = 16.38720414068978
$ws.Range("O11").Value = 25.27882155016518
$ws.Range("B12").Value = 10.22211337607124
$ws.Range("C12").Value = 5.122960783602683
$ws.Range("D12").Value = 10.51136266742621
$ws.Range("F12").Value = 34.29202396621995
$ws.Range("G12").Value = 3.657185409071923
$ws.Range("I12").Value = 23.22575607088029
$ws.Range("J12").Value = 11.37132840823454
$ws.Range("K12").Value = 10.02764174733194
$ws.Range("M12").Value = 16.41667066604305
$ws.Range("O12").Value = 25.26131739806126
$ws.Range("B13").Value = 10.20552199560645
$ws.Range("C13").Value = 5.112428245503043
$ws.Range("D13").Value = 10.50753787855315
$ws.Range("F13").Value = 34.29377359784584
$ws.Range("G13").Value = 3.657279489776339
$ws.Range("I13").Value = 23.23050579514212
$ws.Range("J13").Value = 11.37192102086964
$ws.Range("K13").Value = 10.01616748798559
$ws.Range("M13").Value = 16.41031318738985
$ws.Range("O13").Value = 25.26505339382931
$ws.Range("B14").Value = 10.15124400281495
$ws.Range("C14").Value = 5.077929788037302
$ws.Range("D14").Value = 10.49509750059963
$ws.Range("F14").Value = 34.29967778838589
$ws.Range("G14").Value = 3.657587735655027
$ws.Range("I14").Value = 23.24610308048383
$ws.Range("J14").Value = 11.37389570151012
$ws.Range("K14").Value = 9.978679019482346
$ws.Range("M14").Value = 16.38962323364455
$ws.Range("O14").Value = 25.27736622815426
$ws.Range("B15").Value = 10.11786773797266
$ws.Range("C15").Value = 5.05668402225828
$ws.Range("D15").Value = 10.48750304707528
$ws.Range("F15").Value = 34.30344660368544
$ws.Range("G15").Value = 3.657777643387035
$ws.Range("I15").Value = 23.25573938536413
$ws.Range("J15").Value = 11.37513756139657
$ws.Range("K15").Value = 9.955664521686018
$ws.Range("M15").Value = 16.3769835575492
$ws.Range("O15").Value = 25.28500724238813
$ws.Range("B16").Value = 9.924665884784881
$ws.Range("C16").Value = 4.933198618558602
$ws.Range("D16").Value = 10.44438387776714
$ws.Range("F16").Value = 34.32737801918259
$ws.Range("G16").Value = 3.658882791136452
$ws.Range("I16").Value = 23.31222540572501
$ws.Range("J16").Value = 11.38274959756748
$ws.Range("K16").Value = 9.823020277089874
$ws.Range("M16").Value = 16.30507874022727
$ws.Range("O16").Value = 25.33031159308494
$ws.Range("B17").Value = 9.804469404647463
$ws.Range("C17").Value = 4.855919595893781
$ws.Range("D17").Value = 10.4182992849792
$ws.Range("F17").Value = 34.34414056157677
$ws.Range("G17").Value = 3.659575844080521
$ws.Range("I17").Value = 23.34800531865375
$ws.Range("J17").Value = 11.38786162615139
$ws.Range("K17").Value = 9.741010797088929
$ws.Range("M17").Value = 16.26145411427656
$ws.Range("O17").Value = 25.35945632218456
$ws.Range("B18").Value = 9.734742961502166
$ws.Range("C18").Value = 4.810920124588617
$ws.Range("D18").Value = 10.40343486033609
$ws.Range("F18").Value = 34.35454631810649
$ws.Range("G18").Value = 3.659980025054654
$ws.Range("I18").Value = 23.36899914322385
$ws.Range("J18").Value = 11.39096454353819
$ws.Range("K18").Value = 9.693623161228256
$ws.Range("M18").Value = 16.23654730368337
$ws.Range("O18").Value = 25.37671585260076
$ws.Range("B19").Value = 9.711035213326777
$ws.Range("C19").Value = 4.795590485706938
$ws.Range("D19").Value = 10.39842619256698
$ws.Range("F19").Value = 34.35820075726351
$ws.Range("G19").Value = 3.660117829161462
$ws.Range("I19").Value = 23.37617840496649
$ws.Range("J19").Value = 11.39204307346836
$ws.Range("K19").Value = 9.67754263714364
$ws.Range("M19").Value = 16.22814661734203
$ws.Range("O19").Value = 25.38264480801143
$ws.Range("B20").Value = 9.817326377064003
$ws.Range("C20").Value = 4.864203295904292
$ws.Range("D20").Value = 10.42106175650532
$ws.Range("F20").Value = 34.34227703652656
$ws.Range("G20").Value = 3.659501492761704
$ws.Range("I20").Value = 23.34415361657324
$ws.Range("J20").Value = 11.38730061361542
$ws.Range("K20").Value = 9.749763769473779
$ws.Range("M20").Value = 16.26607902780251
$ws.Range("O20").Value = 25.35630244215912
$ws.Range("B21").Value = 10.16721384912127
$ws.Range("C21").Value = 5.088086721507777
$ws.Range("D21").Value = 10.49874621915763
$ws.Range("F21").Value = 34.29791183746726
$ws.Range("G21").Value = 3.657496967570999
$ws.Range("I21").Value = 23.24150458193348
$ws.Range("J21").Value = 11.37330895685965
$ws.Range("K21").Value = 9.989701155120501
$ws.Range("M21").Value = 16.39569342388329
$ws.Range("O21").Value = 25.2737290055998
$ws.Range("B22").Value = 10.39015126515024
$ws.Range("C22").Value = 5.229309373465687
$ws.Range("D22").Value = 10.55067484271242
$ws.Range("F22").Value = 34.27573316475424
$ws.Range("G22").Value = 3.656236084110493
$ws.Range("I22").Value = 23.17811020090718
$ws.Range("J22").Value = 11.36561097088206
$ws.Range("K22").Value = 10.14424124402307
$ws.Range("M22").Value = 16.48192055921724
$ws.Range("O22").Value = 25.22419404417437
$ws.Range("B23").Value = 10.27171074083618
$ws.Range("C23").Value = 5.154411189918969
$ws.Range("D23").Value = 10.52285750461209
$ws.Range("F23").Value = 34.28694625513358
$ws.Range("G23").Value = 3.656904553937692
$ws.Range("I23").Value = 23.21160686315714
$ws.Range("J23").Value = 11.36958728006072
$ws.Range("K23").Value = 10.06198372890566
$ws.Range("M23").Value = 16.43576711055313
$ws.Range("O23").Value = 25.25022570590317
$ws.Range("B24").Value = 9.811515677828622
$ws.Range("C24").Value = 4.860460010551555
$ws.Range("D24").Value = 10.41981243124939
$ws.Range("F24").Value = 34.34311714171445
$ws.Range("G24").Value = 3.659535089124483
$ws.Range("I24").Value = 23.34589365201464
$ws.Range("J24").Value = 11.38755373669474
$ws.Range("K24").Value = 9.745807292229443
$ws.Range("M24").Value = 16.26398756231456
$ws.Range("O24").Value = 25.35772674201774
$ws.Range("B25").Value = 9.292809672402369
$ws.Range("C25").Value = 4.52250011995731
$ws.Range("D25").Value = 10.31400856365811
$ws.Range("F25").Value = 34.43291937459912
$ws.Range("G25").Value = 3.662585111261054
$ws.Range("I25").Value = 23.50657590295293
$ws.Range("J25").Value = 11.4131509581846
$ws.Range("K25").Value = 9.396651888443222
$ws.Range("M25").Value = 16.08583738842062
$ws.Range("O25").Value = 25.49264887843906
